$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued columns (B, C, D, E) keep their original text type,
# even when the text looks numeric (e.g. "596.38" in column D).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.899.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.635.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.400"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.91"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.104.50"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.725.56"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +12.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.614.67"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.84"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.32%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.76"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.47"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.32"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.53%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "536.31"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.56%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.06"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.79"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.30"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.425"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.72"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.18"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "166.71"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.88"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.42%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.30"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.71%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.23"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +10.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.643"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0252"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.41"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.32%  "
